# Auto-generated from the workbook xml diff: refreshes the crypto table's
# "Price" (column D) and "Volume(1h)" (column E) text cells to the new values.
#
# Every cell in D/E is stored as *text* in the workbook (not a number), even
# when the text looks numeric (e.g. "583.01"). Assigning such a string straight
# to Range.Value makes Excel auto-convert it to a real number, so those cells
# are written with a leading apostrophe (the classic "force text" prefix) and
# then have their Style reset to Normal so the quote-prefix formatting left by
# that trick does not linger on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.994.47'
$ws.Range("E2").Value = '  -0.57%  '
$ws.Range("D3").Value = '2.582.60'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'" + '583.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").Value = "'" + '144.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -2.13%  '
$ws.Range("E9").Value = '  -1.98%  '
$ws.Range("E10").Value = '  -1.90%  '
$ws.Range("E11").Value = '  -0.70%  '
$ws.Range("D12").Value = "'" + '0.351'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.10%  '
$ws.Range("D13").Value = "'" + '27.09'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.71%  '
$ws.Range("D14").Value = '3.045.51'
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("D15").Value = '62.916.68'
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("E16").Value = '  -2.12%  '
$ws.Range("D17").Value = '2.578.01'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").Value = "'" + '11.13'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.24%  '
$ws.Range("D19").Value = "'" + '340.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.93%  '
$ws.Range("E20").Value = '  -1.83%  '
$ws.Range("D21").Value = "'" + '6.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.43%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = "'" + '67.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.98%  '
$ws.Range("E24").Value = '  +7.30%  '
$ws.Range("D25").Value = "'" + '1.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.78%  '
$ws.Range("E26").Value = '  -3.39%  '
$ws.Range("D27").Value = "'" + '8.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.87%  '
$ws.Range("D28").Value = "'" + '0.997'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.21%  '
$ws.Range("D29").Value = "'" + '8.27'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.40%  '
$ws.Range("E30").Value = '  -2.91%  '
$ws.Range("D31").Value = "'" + '459.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").Value = '0.0₃0801'
$ws.Range("E32").Value = '  -3.34%  '
$ws.Range("D33").Value = "'" + '1.66'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("D34").Value = "'" + '176.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -2.37%  '
$ws.Range("D37").Value = "'" + '18.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.09%  '
$ws.Range("D38").Value = "'" + '4.49'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.82%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D41").Value = "'" + '159.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.66%  '
$ws.Range("D42").Value = "'" + '40.07'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.91%  '
$ws.Range("D43").Value = "'" + '3.70'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.10%  '
$ws.Range("D44").Value = "'" + '21.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("E45").Value = '  +3.14%  '
$ws.Range("D46").Value = "'" + '0.0539'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.86%  '
$ws.Range("D47").Value = "'" + '0.0962'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.20%  '
$ws.Range("E48").Value = '  -1.82%  '
$ws.Range("D49").Value = "'" + '18.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.72%  '
$ws.Range("D50").Value = "'" + '11.41'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("E51").Value = '  -4.47%  '
